# LOB1262.docx restructuring script
#
# The commit reshuffles several paragraphs of text: the same literal runs
# of text end up in different paragraphs, several of them forming a
# rotation (A's old text becomes B's new text, B's old text becomes C's,
# ...). A naive sequence of Find&Replace calls would make a later search
# string ambiguous once an earlier step has produced a duplicate copy of
# it elsewhere in the document. To make every Find.Execute unconditionally
# safe we do this in two phases:
#
#   Phase 1: replace every distinct "old" text (each is unique in the
#            pristine document) with a unique placeholder token.
#   Phase 2: replace every placeholder token (now unique, since we just
#            put it there) with its final text.
#
# Paragraph count / order / styles / run formatting (bold, italic) do not
# change anywhere, so plain literal Find&Replace (wildcards off) safely
# keeps existing run boundaries, <w:br/> separators and rPr intact.

$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $range = $d.Content
    $ok = $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $ok) {
        Write-Host ("NOT FOUND: " + $old)
    }
}

# ----------------------------------------------------------------------
# The eleven distinct text values that move around, keyed by a short
# mnemonic, exactly as they appear in the source document today.
# ----------------------------------------------------------------------

$V_objetivo_pt  = "Propiciar ao discente conhecimento dos fundamentos da Educação Ambiental utilizando como base os problemas ambientais da atualidade. Desenvolver atividades práticas integradas à região. Orientar o desenvolvimento de projetos relacionados à Gestão e Educação Ambiental."
$V_objetivo_en  = "Provide students with knowledge of the fundamentals of Environmental Education using as basis the current environmental problems. To develop practical activities integrated to the region. Guide the development of projects related to Environmental Education and Management"
$V_docente1     = "9146830 - Danúbia Caporusso Bargos"
$V_docente2     = "5817650 - Érica Leonor Romão"
$V_resumo_pt    = "Considerações gerais sobre a problemática ambiental. Evolução das questões ambientais no Brasil e no mundo. Educação e Gestão Ambiental. Elaboração e acompanhamento de projetos de educação ambiental."
$V_resumo_en    = "General considerations on environmental problem. Evolution of environmental questions in Brazil and in the world. Education and Environmental Management. Development and monitoring of environmental education projects."
$V_programa_pt  = "Sociedade, natureza e desenvolvimento. A relação degradação ambiental-qualidade de vida. Meio ambiente e cidadania. Percepção e Interpretação ambiental. Meio ambiente e representação social. Histórico da educação ambiental e conceitos de meio ambiente; Conceitos, princípios e pensamentos norteadores da Educação Ambiental. A questão ambiental e as conferências mundiais de meio ambiente. O movimento ambientalista e o histórico da EA no Brasil e no mundo; A Agenda 21 e educação ambiental. A política nacional de educação ambiental (pnea) e legislação correlata: A abordagem interdisciplinar da educação ambiental; Educação como instrumento de Gestão Ambiental. Educação ambiental nas empresas e o Sistema de Gestão Ambiental. Projetos, reflexões e práticas da Educação Ambiental. Análise e vivências de experiências práticas de educação ambiental em diferentes contextos. Metodologia de projetos, oficinas e capacitação em educação ambiental."
$V_metodo_val   = "Avaliação baseada em provas, exercícios, projetos, seminários e outras formas de avaliação, sendo a nota final correspondente a média ponderada das notas atribuídas às avaliações aplicadas"
$V_criterio_val = "Nota Final: NF ≥ 5,0"
$V_norma_val    = "Provas e/ou exercícios dirigidos"
$V_biblio       = "CARVALHO, I. C. M.; Educação Ambiental e formação do sujeito ecológico. São Paulo: Cortez, 2006." + [char]11 + `
                   "CINQUETTI, H. C. S.; LOGAREZZI, A. (Org.). Consumo e Resíduo - Fundamentos para o trabalho educativo. 1 ed. São Carlos: EdUFSCar, 2006, v. 1." + [char]11 + `
                   "DIAS, G. F. Dinâmica e instrumentação para educação ambiental. 1. ed. São Paulo: Gaia, 2010. v. 1. 216p." + [char]11 + `
                   "DIAS, G. F. Educação e Gestão Ambiental. 1. ed. São Paulo: Editora Gaia Ltda, 2006. v. 1. 118p." + [char]11 + `
                   "DIAS, G. F. Educação Ambiental: princípios e práticas. 6a ed. São Paulo: Gaia, 2000." + [char]11 + `
                   "GUIMARÃES, M. (org.) Caminhos da educação ambiental: da forma à ação. Campinas, SP: Papirus, 2006." + [char]11 + `
                   "JACOBI, Pedro Roberto, MONTEIRO,F. M ; FERNANDES, M. L. B. . Educação e Sustentabilidade- caminhos e práticas para uma educação transformadora. São Paulo: Evoluir Cultural, 2009. v. 01. 108p." + [char]11 + `
                   "JACOBI, Pedro Roberto OLIVEIRA, F. C. J. F. (Org.). Educação, Meio Ambiente e Cidadania - reflexões e experiências. São Paulo: SMA/CEAM, 1998. 121p " + [char]11 + `
                   "LOUREIRO, C. F. B. Trajetória e fundamentos da educação ambiental. 4. ed. São Paulo: Cortez editora, 2012. 165p" + [char]11 + `
                   "PHILIPPI JR., A & PELICIONI, M. C. F. (Eds). 2005. Educação ambiental e sustentabilidade. Barueri SP: Manole. 878p. (Coleção Ambiental, 3)."

# ----------------------------------------------------------------------
# Phase 1: stash every old value behind a unique placeholder token.
# ----------------------------------------------------------------------

Replace-Text $V_objetivo_pt  "@@PH_OBJETIVO_PT@@"
Replace-Text $V_objetivo_en  "@@PH_OBJETIVO_EN@@"
Replace-Text $V_docente1     "@@PH_DOCENTE1@@"
Replace-Text $V_docente2     "@@PH_DOCENTE2@@"
Replace-Text $V_resumo_pt    "@@PH_RESUMO_PT@@"
Replace-Text $V_resumo_en    "@@PH_RESUMO_EN@@"
Replace-Text $V_programa_pt  "@@PH_PROGRAMA_PT@@"
Replace-Text $V_metodo_val   "@@PH_METODO@@"
Replace-Text $V_criterio_val "@@PH_CRITERIO@@"
Replace-Text $V_norma_val    "@@PH_NORMA@@"
Replace-Text $V_biblio       "@@PH_BIBLIO@@"

# ----------------------------------------------------------------------
# Phase 2: each placeholder becomes the real final text at that spot.
# ----------------------------------------------------------------------

Replace-Text "@@PH_OBJETIVO_PT@@" $V_resumo_pt     # P6  : Objetivos (PT)        <- Programa resumido (PT)
Replace-Text "@@PH_OBJETIVO_EN@@" $V_resumo_en     # P7  : Objetivos (EN)        <- Programa resumido (EN)
Replace-Text "@@PH_DOCENTE1@@"    $V_objetivo_pt   # P9  run1 : Docente list     <- old Objetivos (PT)
Replace-Text "@@PH_DOCENTE2@@"    $V_programa_pt   # P9  run2 : Docente list     <- old Programa (PT)
Replace-Text "@@PH_RESUMO_PT@@"   $V_metodo_val    # P11 : Programa resumido (PT)<- old Método value
Replace-Text "@@PH_RESUMO_EN@@"   $V_objetivo_en   # P12 : Programa resumido (EN)<- old Objetivos (EN)
Replace-Text "@@PH_PROGRAMA_PT@@" $V_criterio_val  # P14 : Programa (PT)         <- old Critério value
Replace-Text "@@PH_METODO@@"      $V_norma_val     # P17 run2 : Método value     <- old Norma value
Replace-Text "@@PH_CRITERIO@@"    $V_biblio        # P17 run4 : Critério value   <- old Bibliografia body
Replace-Text "@@PH_NORMA@@"       $V_docente1      # P17 run6 : Norma value      <- old Docente1
Replace-Text "@@PH_BIBLIO@@"      $V_docente2      # P19 : Bibliografia body     <- old Docente2

Write-Host "done"
